$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 306, shifting existing rows 306-338 down to 307-339.
$ws.Rows.Item(306).Insert()

# Populate the newly inserted row 306 with the new record's data.
$ws.Cells.Item(306, 1).Value = 5
$ws.Cells.Item(306, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(306, 3).Value = "Maule"
$ws.Cells.Item(306, 4).Value = 45132
$ws.Cells.Item(306, 5).Value = 7
$ws.Cells.Item(306, 6).Value = 100112017
$ws.Cells.Item(306, 7).Value = "Apio"
$ws.Cells.Item(306, 8).Value = "Americana (o)"
$ws.Cells.Item(306, 9).Value = "Primera"
$ws.Cells.Item(306, 10).Value = 700
$ws.Cells.Item(306, 11).Value = 5500
$ws.Cells.Item(306, 12).Value = 5500
$ws.Cells.Item(306, 13).Value = 5500
$ws.Cells.Item(306, 14).Value = "$/docena de matas"
$ws.Cells.Item(306, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(306, 16).Value = 917
$ws.Cells.Item(306, 17).Value = 6
$ws.Cells.Item(306, 18).Value = "Hortaliza"

Write-Output "done"
